# "Doing Updates for Financials"
# Updates the yearly financial figures on the KZIA sheet (columns D:J hold the
# 7 reported fiscal years). Only the numeric cell values changed; formatting,
# layout, and labels are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KZIA")

# --- Income Statement -------------------------------------------------

# Research Development
$ws.Range("D12").Value = 6900
$ws.Range("E12").Value = 7900
$ws.Range("F12").Value = 7000
$ws.Range("G12").Value = 4200

# Non Recurring
$ws.Range("D14").Value = -5000

# Total Operating Expenses
$ws.Range("D17").Value = 3900
$ws.Range("E17").Value = 13400
$ws.Range("F17").Value = 11500
$ws.Range("G17").Value = 6900
$ws.Range("H17").Value = 4800
$ws.Range("I17").Value = 2200

# Operating Income or Loss
$ws.Range("E18").Value = -13200
$ws.Range("F18").Value = -11200
$ws.Range("G18").Value = -6900
$ws.Range("H18").Value = -4700
$ws.Range("J18").Value = -1000

# Total Other Income/Expenses Net
$ws.Range("E20").Value = 5500
$ws.Range("F20").Value = 2600
$ws.Range("I20").Value = 300

# Earnings Before Interest And Taxes
$ws.Range("D21").Value = -3400
$ws.Range("E21").Value = -6700
$ws.Range("F21").Value = -8100
$ws.Range("G21").Value = -4800
$ws.Range("H21").Value = -4600
$ws.Range("I21").Value = -800
$ws.Range("J21").Value = -1000

# Interest Expense
$ws.Range("H22").Value = 300

# Income Before Tax
$ws.Range("D23").Value = -4500
$ws.Range("E23").Value = -7700
$ws.Range("F23").Value = -8600
$ws.Range("G23").Value = -5200
$ws.Range("H23").Value = -5400
$ws.Range("J23").Value = -1000

# Income After Tax
$ws.Range("D26").Value = -4300
$ws.Range("E26").Value = -7500
$ws.Range("F26").Value = -8600
$ws.Range("G26").Value = -5200
$ws.Range("H26").Value = -5400
$ws.Range("J26").Value = -1000

# Net Income From Continuing Ops
$ws.Range("D27").Value = -4300
$ws.Range("E27").Value = -7500
$ws.Range("F27").Value = -8500
$ws.Range("G27").Value = -5100
$ws.Range("H27").Value = -5300
$ws.Range("J27").Value = 800

# Other Items
$ws.Range("E32").Value = -5500
$ws.Range("F32").Value = -2600
$ws.Range("I32").Value = -300

# Net Income
$ws.Range("D33").Value = -4300
$ws.Range("E33").Value = -7500
$ws.Range("F33").Value = -8500
$ws.Range("G33").Value = -5100
$ws.Range("H33").Value = -5300

# Net Income Applicable To Common Shares
$ws.Range("D35").Value = -4300
$ws.Range("E35").Value = -7500
$ws.Range("F35").Value = -8500
$ws.Range("G35").Value = -5100
$ws.Range("H35").Value = -5300

# --- Balance Sheet ------------------------------------------------------

# Cash And Cash Equivalents
$ws.Range("D41").Value = 2100
$ws.Range("E41").Value = 10200
$ws.Range("F41").Value = 23700
$ws.Range("G41").Value = 31400
$ws.Range("J41").Value = 4500

# Short Term Investments (D42 switches from "NA" to a real figure)
$ws.Range("D42").Value = 2100
$ws.Range("J42").Value = 1200

# Net Receivables
$ws.Range("E43").Value = 3000

# Other Current Assets
$ws.Range("D45").Value = 500
$ws.Range("E45").Value = 500

# Total Current Assets
$ws.Range("D46").Value = 6600
$ws.Range("E46").Value = 13800
$ws.Range("F46").Value = 24100
$ws.Range("G46").Value = 31600
$ws.Range("I46").Value = 2200
$ws.Range("J46").Value = 6300

# Property Plant and Equipment
$ws.Range("E48").Value = 300

# Goodwill
$ws.Range("D49").Value = 10300
$ws.Range("E49").Value = 11300

# Total Assets
$ws.Range("D54").Value = 19900
$ws.Range("E54").Value = 25400
$ws.Range("F54").Value = 25100
$ws.Range("G54").Value = 32600
$ws.Range("H54").Value = 3300
$ws.Range("I54").Value = 4100
$ws.Range("J54").Value = 6400

# Accounts Payable
$ws.Range("D57").Value = 1000
$ws.Range("E57").Value = 1300
$ws.Range("G57").Value = 1100

# Short/Current Long Term Debt
$ws.Range("H58").Value = 1900

# Other Current Liabilities
$ws.Range("D59").Value = 1800
$ws.Range("J59").Value = 2100

# Total Current Liabilities
$ws.Range("E60").Value = 3800
$ws.Range("H60").Value = 2300
$ws.Range("J60").Value = 2700

# Other Liabilities
$ws.Range("D62").Value = 3600
$ws.Range("E62").Value = 3700

# Total Liabilities
$ws.Range("D66").Value = 6300
$ws.Range("E66").Value = 7500
$ws.Range("F66").Value = 1100
$ws.Range("G66").Value = 1000
$ws.Range("H66").Value = 2200
$ws.Range("J66").Value = 3900

# Retained Earnings
$ws.Range("D72").Value = -8800
$ws.Range("E72").Value = -119600
$ws.Range("F72").Value = -112600
$ws.Range("G72").Value = -104300
$ws.Range("H72").Value = -100000
$ws.Range("I72").Value = -94500
$ws.Range("J72").Value = -138400

# Total Stockholder Equity
$ws.Range("D76").Value = 13600
$ws.Range("E76").Value = 17900
$ws.Range("F76").Value = 24000
$ws.Range("G76").Value = 31600

# --- Cash Flow Statement -------------------------------------------------

# Net Income
$ws.Range("D81").Value = -4300
$ws.Range("E81").Value = -7500
$ws.Range("F81").Value = -8500
$ws.Range("G81").Value = -5100
$ws.Range("H81").Value = -5300

# Total Cash Flow From Operating Activities
$ws.Range("D89").Value = -6100
$ws.Range("E89").Value = -8100
$ws.Range("F89").Value = -8500
$ws.Range("G89").Value = -4100
$ws.Range("H89").Value = -4000
$ws.Range("I89").Value = -6200
$ws.Range("J89").Value = -5200

# Total Cash Flows From Investing Activities
$ws.Range("E94").Value = -5000
$ws.Range("J94").Value = 5700

# Total Cash Flows From Financing Activities
$ws.Range("G100").Value = 33500
$ws.Range("H100").Value = 3900
$ws.Range("I100").Value = 2100
$ws.Range("J100").Value = 1200

# Change In Cash and Cash Equivalents
$ws.Range("D102").Value = -6000
$ws.Range("E102").Value = -13400
$ws.Range("F102").Value = -7700
$ws.Range("G102").Value = 29600
$ws.Range("I102").Value = -4000
$ws.Range("J102").Value = 1600
